# daily auto push: 2026-02-07 02:49 UTC
# Insert a new data row for 2026/02/07 08:00 (time-slot 8) just before the
# 2026/12/29 block, shifting the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 784; everything currently at/after row 784 (starting
# with 2026/12/29) shifts down to make room.
$ws.Rows.Item(784).Insert()

# Column A holds a date-like text string ("yyyy/mm/dd"), not a real date
# value. Force the cell to text formatting first so Excel doesn't
# auto-convert the literal into a date serial number, then drop the format
# back to the sheet's default (General) once the text value is in place.
$ws.Range("A784").NumberFormat = "@"
$ws.Range("A784").Value = "2026/02/07"
$ws.Range("A784").Style = "Normal"

$ws.Range("B784").Value = "土"
$ws.Range("C784").Value = 8
$ws.Range("D784").Value = 201
